$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet
$ws.Name = "Pizzas a la semana"

# Remove the chart/drawing objects from the sheet
foreach ($chartObj in $ws.ChartObjects()) {
    $chartObj.Delete()
}
